# Improved French sentiment analysis and fixed UI issues
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix "Secteur d'activité" (column Y) labels: "Industrie Automobile" -> "Industrie " ---
$industrieRows = @(12, 16, 17, 18, 19, 20)
foreach ($r in $industrieRows) {
    $ws.Cells.Item($r, 25).Value = "Industrie "
}

# --- Update "Advancement rate" (column W) values ---
$ws.Range("W2").Value = 33
$ws.Range("W3").Value = 12
$ws.Range("W4").Value = 77
$ws.Range("W5").Value = 10
$ws.Range("W6").Value = 9
$ws.Range("W7").Value = 44
$ws.Range("W8").Value = 22
$ws.Range("W9").Value = 35
$ws.Range("W10").Value = 88
$ws.Range("W11").Value = 78
$ws.Range("W12").Value = 87
$ws.Range("W13").Value = 54
$ws.Range("W14").Value = 23
$ws.Range("W15").Value = 4
$ws.Range("W16").Value = 98
$ws.Range("W17").Value = 66
$ws.Range("W18").Value = 61
$ws.Range("W19").Value = 28
$ws.Range("W20").Value = 36

# --- Fix UI: scroll position + active selection ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 19
$win.ScrollRow = 1
$ws.Range("W20").Select()
